$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.04879550243866605
$ws.Range("D2").Value = 0.005954525361048833
$ws.Range("E2").Value = 1.645051762268533

# Row 3
$ws.Range("C3").Value = 0.1884369884805298
$ws.Range("D3").Value = 0.1427852673396272
$ws.Range("E3").Value = 1.046585486477964

# Row 4
$ws.Range("C4").Value = 0.1177273822883147
$ws.Range("D4").Value = 0.1097642259820403
$ws.Range("E4").Value = 0.7172408094163223

# Row 5
$ws.Range("B5").Value = 0.1875320075471793
$ws.Range("C5").Value = 0.1197300911685313
$ws.Range("D5").Value = 0.1143798086650707
$ws.Range("E5").Value = 0.7058071884429076

# Row 6
$ws.Range("B6").Value = 0.152155804319207
$ws.Range("C6").Value = 0.2281010324356234
$ws.Range("D6").Value = 0.2338861163958089
$ws.Range("E6").Value = 0.8085175612373885

# Row 7
$ws.Range("B7").Value = 0.06031218813333689
$ws.Range("C7").Value = 0.1090459691992094
$ws.Range("D7").Value = 0.1481662248229568
$ws.Range("E7").Value = 0.472752607977338
